$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 29
$ws.Range("H29").Value = 2877.9443
# Row 38
$ws.Range("H38").Value = 561
$ws.Range("I38").Value = 184.85715
$ws.Range("J38").Value = 999.8333
$ws.Range("K38").Value = 554.5714499999999
$ws.Range("L38").Value = 2999.4999
$ws.Range("M38").Value = -182.5714499999999
$ws.Range("N38").Value = -3743.4999
# Row 47
$ws.Range("H47").Value = 15000
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
# Row 58
$ws.Range("H58").Value = 2262.5
$ws.Range("J58").Value = 2450
$ws.Range("L58").Value = 7350
$ws.Range("N58").Value = -7650
# Row 86
$ws.Range("H86").Value = 39279.63
$ws.Range("I86").Value = 43472.957
$ws.Range("J86").Value = 5733
$ws.Range("K86").Value = 43472.957
$ws.Range("L86").Value = 5733
$ws.Range("M86").Value = -42349.957
$ws.Range("N86").Value = -7979
# Row 87
$ws.Range("H87").Value = 28830.8
$ws.Range("J87").Value = 28830.8
$ws.Range("L87").Value = 28830.8
$ws.Range("N87").Value = -31326.8
# Row 89
$ws.Range("H89").Value = 39279.63
$ws.Range("I89").Value = 43472.957
$ws.Range("J89").Value = 5733
$ws.Range("K89").Value = 217364.785
$ws.Range("L89").Value = 28665
$ws.Range("M89").Value = -211748.785
$ws.Range("N89").Value = -39897
# Row 90
$ws.Range("H90").Value = 28830.8
$ws.Range("J90").Value = 28830.8
$ws.Range("L90").Value = 86492.39999999999
$ws.Range("N90").Value = -98972.39999999999
# Row 121
$ws.Range("H121").Value = 1756.8
$ws.Range("I121").Value = 2450
$ws.Range("J121").Value = 1294.6666
$ws.Range("K121").Value = 7350
$ws.Range("L121").Value = 3883.9998
$ws.Range("M121").Value = -5603
$ws.Range("N121").Value = -7377.9998
# Row 132
$ws.Range("H132").Value = 1987.9412
$ws.Range("I132").Value = 2062.1875
$ws.Range("K132").Value = 6186.5625
$ws.Range("M132").Value = -3656.5625
# Row 137
$ws.Range("H137").Value = 507443.94
$ws.Range("I137").Value = 2521.2903
$ws.Range("J137").Value = 954661.1
$ws.Range("K137").Value = 7563.8709
$ws.Range("L137").Value = 2863983.3
$ws.Range("M137").Value = -5013.8709
$ws.Range("N137").Value = -2869083.3

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 19464.05
$ws.Range("I32").Value = 20425.518
$ws.Range("J32").Value = 6003.5
$ws.Range("K32").Value = 20425.518
$ws.Range("L32").Value = 6003.5
$ws.Range("M32").Value = -20138.518
$ws.Range("N32").Value = -6577.5
# Row 60
$ws.Range("H60").Value = 29000
$ws.Range("I60").Value = 29000
$ws.Range("K60").Value = 29000
$ws.Range("M60").Value = -28267
# Row 61
$ws.Range("H61").Value = 6361.311
$ws.Range("I61").Value = 2790.5667
$ws.Range("K61").Value = 2790.5667
$ws.Range("M61").Value = -2578.5667
# Row 132
$ws.Range("H132").Value = 3371.7144
$ws.Range("I132").Value = 2383.5
$ws.Range("J132").Value = 4112.875
$ws.Range("K132").Value = 7150.5
$ws.Range("L132").Value = 12338.625
$ws.Range("M132").Value = -4620.5
$ws.Range("N132").Value = -17398.625
# Row 136
$ws.Range("H136").Value = 6361.311
$ws.Range("I136").Value = 2790.5667
$ws.Range("K136").Value = 8371.7001
$ws.Range("M136").Value = -5821.7001

$ws = $wb.Worksheets.Item("BSM")
# Row 108
$ws.Range("H108").Value = 60684
$ws.Range("J108").Value = 60684
$ws.Range("L108").Value = 60684
$ws.Range("N108").Value = -68364
# Row 111
$ws.Range("H111").Value = 36599.332
$ws.Range("J111").Value = 36599.332
$ws.Range("L111").Value = 36599.332
$ws.Range("N111").Value = -44779.332

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 563116.4
$ws.Range("I31").Value = 11843.056
$ws.Range("J31").Value = 778832
$ws.Range("K31").Value = 11843.056
$ws.Range("L31").Value = 778832
$ws.Range("M31").Value = -11548.056
$ws.Range("N31").Value = -779422
# Row 34
$ws.Range("H34").Value = 563116.4
$ws.Range("I34").Value = 11843.056
$ws.Range("J34").Value = 778832
$ws.Range("K34").Value = 11843.056
$ws.Range("L34").Value = 778832
$ws.Range("M34").Value = -11641.056
$ws.Range("N34").Value = -779236
# Row 141
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 117362.06
$ws.Range("I68").Value = 217484.44
$ws.Range("J68").Value = 4180.2393
$ws.Range("K68").Value = 652453.3200000001
$ws.Range("L68").Value = 12540.7179
$ws.Range("M68").Value = -651642.3200000001
$ws.Range("N68").Value = -14162.7179
# Row 71
$ws.Range("H71").Value = 117362.06
$ws.Range("I71").Value = 217484.44
$ws.Range("J71").Value = 4180.2393
$ws.Range("K71").Value = 1957359.96
$ws.Range("L71").Value = 37622.1537
$ws.Range("M71").Value = -1953303.96
$ws.Range("N71").Value = -45734.1537
# Row 107
$ws.Range("H107").Value = 1417.8723
$ws.Range("I107").Value = 767.3077
$ws.Range("J107").Value = 1666.6177
$ws.Range("K107").Value = 2301.9231
$ws.Range("L107").Value = 4999.8531
$ws.Range("M107").Value = -381.9231
$ws.Range("N107").Value = -8839.8531
# Row 132
$ws.Range("H132").Value = 5980
$ws.Range("I132").Value = 10000
$ws.Range("J132").Value = 4975
$ws.Range("K132").Value = 90000
$ws.Range("L132").Value = 44775
$ws.Range("M132").Value = -87470
$ws.Range("N132").Value = -49835

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5605.482
$ws.Range("I70").Value = 5427.3145
$ws.Range("J70").Value = 5902.4287
$ws.Range("K70").Value = 5427.3145
$ws.Range("L70").Value = 5902.4287
$ws.Range("M70").Value = -5157.3145
$ws.Range("N70").Value = -6442.4287
# Row 73
$ws.Range("H73").Value = 5605.482
$ws.Range("I73").Value = 5427.3145
$ws.Range("J73").Value = 5902.4287
$ws.Range("K73").Value = 5427.3145
$ws.Range("L73").Value = 5902.4287
$ws.Range("M73").Value = -4491.3145
$ws.Range("N73").Value = -7774.4287

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 2250
$ws.Range("I93").Value = 1000
$ws.Range("J93").Value = 3500
$ws.Range("K93").Value = 1000
$ws.Range("L93").Value = 3500
$ws.Range("M93").Value = 248
$ws.Range("N93").Value = -5996
# Row 136
$ws.Range("H136").Value = 5269.3335
$ws.Range("I136").Value = 2227.7222
$ws.Range("J136").Value = 7876.4287
$ws.Range("K136").Value = 6683.1666
$ws.Range("L136").Value = 23629.2861
$ws.Range("M136").Value = -4133.1666
$ws.Range("N136").Value = -28729.2861
# Row 141
$ws.Range("H141").Value = 59357.5
$ws.Range("J141").Value = 59357.5
$ws.Range("L141").Value = 59357.5
$ws.Range("N141").Value = -69717.5

$ws = $wb.Worksheets.Item("WVR")
# Row 6
$ws.Range("H6").Value = 3335.3333
$ws.Range("J6").Value = 3335.3333
$ws.Range("L6").Value = 3335.3333
$ws.Range("N6").Value = -3565.3333
# Row 12
$ws.Range("H12").Value = 10000000
$ws.Range("I12").Value = 10000000
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 10000000
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -9999858
$ws.Range("N12").ClearContents()
# Row 32
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
# Row 62
$ws.Range("H62").Value = 3988.889
$ws.Range("J62").Value = 3985.7144
$ws.Range("L62").Value = 3985.7144
$ws.Range("N62").Value = -5233.7144
# Row 65
$ws.Range("H65").Value = 3988.889
$ws.Range("J65").Value = 3985.7144
$ws.Range("L65").Value = 19928.572
$ws.Range("N65").Value = -26168.572
# Row 69
$ws.Range("H69").Value = 36567.375
$ws.Range("J69").Value = 36567.375
$ws.Range("L69").Value = 36567.375
$ws.Range("N69").Value = -38065.375
# Row 72
$ws.Range("H72").Value = 36567.375
$ws.Range("J72").Value = 36567.375
$ws.Range("L72").Value = 109702.125
$ws.Range("N72").Value = -117190.125
# Row 92
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

